$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) After the "https://www.apache.org/licenses/LICENSE-2.0" paragraph,
#    insert: one blank paragraph, then the Apache.txt link paragraph,
#    then the Apache.pdf link paragraph (all in the "Roboto Mono" run
#    style used by the license-link paragraphs in this document).
# -----------------------------------------------------------------------
$findLicense = $d.Content
$foundLicense = $findLicense.Find.Execute("licenses/LICENSE-2.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundLicense) {
    $licenseParaIndex = $findLicense.Paragraphs.Item(1).Index

    # Blank paragraph right after the LICENSE-2.0 link.
    $licensePara = $d.Paragraphs.Item($licenseParaIndex)
    $licensePara.Range.InsertParagraphAfter()

    # Paragraph with the Apache.txt link.
    $blankPara = $d.Paragraphs.Item($licenseParaIndex + 1)
    $blankPara.Range.InsertParagraphAfter()
    $txtPara = $d.Paragraphs.Item($licenseParaIndex + 2)
    $txtPara.Range.Text = "https://www.nicolesharp.net/licenses/Apache.txt"

    # Paragraph with the Apache.pdf link.
    $txtPara2 = $d.Paragraphs.Item($licenseParaIndex + 2)
    $txtPara2.Range.InsertParagraphAfter()
    $pdfPara = $d.Paragraphs.Item($licenseParaIndex + 3)
    $pdfPara.Range.Text = "https://www.nicolesharp.net/licenses/Apache.pdf"
}

# -----------------------------------------------------------------------
# 2) Before the "https://www.creativecommons.org/licenses/by-sa/4.0/"
#    paragraph, insert a new paragraph with the CC_BY-SA.htm link
#    (again in the "Roboto Mono" run style already used by that
#    paragraph, which InsertParagraphBefore picks up automatically).
# -----------------------------------------------------------------------
$findCC = $d.Content
$foundCC = $findCC.Find.Execute("creativecommons.org/licenses/by-sa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundCC) {
    $ccParaIndex = $findCC.Paragraphs.Item(1).Index

    $ccPara = $d.Paragraphs.Item($ccParaIndex)
    $ccPara.Range.InsertParagraphBefore()
    $newCCPara = $d.Paragraphs.Item($ccParaIndex)
    $newCCPara.Range.Text = "https://www.nicolesharp.net/licenses/CC_BY-SA.htm"
}
